$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27; existing rows 27..64 shift down to 28..65.
$ws.Rows.Item(27).Insert()

# Populate the new row 27 with the new weekly record.
$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Vega Monumental Concepción"
$ws.Range("C27").Value = "Bíobío"
$ws.Range("D27").Value = 45175
$ws.Range("E27").Value = 8
$ws.Range("F27").Value = 100112026
$ws.Range("G27").Value = "Haba"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 14000
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = 14000
$ws.Range("N27").Value = "$/saco 25 kilos"
$ws.Range("O27").Value = "Región de Coquimbo"
$ws.Range("P27").Value = 560
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = "Hortaliza"
